# Auto-generated edit script: apply updated Leve-profit figures per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 71770.414
$ws.Range("J3").Value = 71770.414
$ws.Range("L3").Value = 71770.414
$ws.Range("N3").Value = -71998.414
$ws.Range("H5").Value = 6671711
$ws.Range("I5").Value = 4006053.2
$ws.Range("J5").Value = 20000000
$ws.Range("K5").Value = 4006053.2
$ws.Range("L5").Value = 20000000
$ws.Range("M5").Value = -4005938.2
$ws.Range("N5").Value = -20000230
$ws.Range("H12").Value = 350.22223
$ws.Range("I12").Value = 385.25
$ws.Range("K12").Value = 385.25
$ws.Range("M12").Value = -215.25
$ws.Range("H53").Value = 336.05884
$ws.Range("J53").Value = 163.5
$ws.Range("L53").Value = 163.5
$ws.Range("N53").Value = -1437.5
$ws.Range("H62").Value = 3998.25
$ws.Range("I62").Value = 3996.5
$ws.Range("K62").Value = 3996.5
$ws.Range("M62").Value = -3372.5
$ws.Range("H65").Value = 3998.25
$ws.Range("I65").Value = 3996.5
$ws.Range("K65").Value = 19982.5
$ws.Range("M65").Value = -16862.5
$ws.Range("H98").Value = 5243.636
$ws.Range("I98").Value = 5243.636
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 5243.636
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -3745.636
$ws.Range("N98").ClearContents()
$ws.Range("H99").Value = 772.1667
$ws.Range("I99").Value = 766
$ws.Range("J99").Value = 784.5
$ws.Range("K99").Value = 2298
$ws.Range("L99").Value = 2353.5
$ws.Range("M99").Value = -800
$ws.Range("N99").Value = -5349.5
$ws.Range("H102").Value = 71770.414
$ws.Range("J102").Value = 71770.414
$ws.Range("L102").Value = 71770.414
$ws.Range("N102").Value = -78260.414
$ws.Range("H122").Value = 5243.636
$ws.Range("I122").Value = 5243.636
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15730.908
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -13280.908
$ws.Range("N122").ClearContents()
$ws.Range("H137").Value = 1642.1875
$ws.Range("I137").Value = 1521.25
$ws.Range("K137").Value = 4563.75
$ws.Range("M137").Value = -2013.75
$ws.Range("H138").Value = 5281.447
$ws.Range("J138").Value = 5606.6665
$ws.Range("L138").Value = 16819.9995
$ws.Range("N138").Value = -27099.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3045.3284
$ws.Range("I32").Value = 2199.0164
$ws.Range("K32").Value = 2199.0164
$ws.Range("M32").Value = -1912.0164
$ws.Range("H45").Value = 965.3855600000001
$ws.Range("I45").Value = 913.6582
$ws.Range("J45").Value = 1987
$ws.Range("K45").Value = 913.6582
$ws.Range("L45").Value = 1987
$ws.Range("M45").Value = -536.6582
$ws.Range("N45").Value = -2741
$ws.Range("H101").Value = 91000
$ws.Range("J101").Value = 91000
$ws.Range("L101").Value = 91000
$ws.Range("N101").Value = -97490

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 61369
$ws.Range("J60").Value = 61369
$ws.Range("L60").Value = 61369
$ws.Range("N60").Value = -62567

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 705.84
$ws.Range("I7").Value = 339.66666
$ws.Range("J7").Value = 911.8125
$ws.Range("K7").Value = 339.66666
$ws.Range("L7").Value = 911.8125
$ws.Range("M7").Value = -226.66666
$ws.Range("N7").Value = -1137.8125
$ws.Range("H22").Value = 2023.5454
$ws.Range("I22").Value = 1753
$ws.Range("K22").Value = 1753
$ws.Range("M22").Value = -1403
$ws.Range("H31").Value = 2457.8823
$ws.Range("I31").Value = 2323.1
$ws.Range("K31").Value = 2323.1
$ws.Range("M31").Value = -2028.1
$ws.Range("H34").Value = 2457.8823
$ws.Range("I34").Value = 2323.1
$ws.Range("K34").Value = 2323.1
$ws.Range("M34").Value = -2121.1
$ws.Range("H86").Value = 5623
$ws.Range("I86").Value = 5333.3335
$ws.Range("K86").Value = 5333.3335
$ws.Range("M86").Value = -4210.3335
$ws.Range("H89").Value = 5623
$ws.Range("I89").Value = 5333.3335
$ws.Range("K89").Value = 26666.6675
$ws.Range("M89").Value = -21050.6675
$ws.Range("H99").Value = 23446.79
$ws.Range("I99").Value = 23765.777
$ws.Range("K99").Value = 23765.777
$ws.Range("M99").Value = -22267.777
$ws.Range("H122").Value = 4399.857
$ws.Range("I122").Value = 4399.857
$ws.Range("K122").Value = 13199.571
$ws.Range("M122").Value = -10749.571
$ws.Range("H126").Value = 23446.79
$ws.Range("I126").Value = 23765.777
$ws.Range("K126").Value = 71297.33099999999
$ws.Range("M126").Value = -68827.33099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4048.1458
$ws.Range("J68").Value = 4180.7827
$ws.Range("L68").Value = 12542.3481
$ws.Range("N68").Value = -14164.3481
$ws.Range("H71").Value = 4048.1458
$ws.Range("J71").Value = 4180.7827
$ws.Range("L71").Value = 37627.04429999999
$ws.Range("N71").Value = -45739.04429999999
$ws.Range("H107").Value = 1065.4
$ws.Range("I107").Value = 898.25
$ws.Range("J107").Value = 1206.1578
$ws.Range("K107").Value = 2694.75
$ws.Range("L107").Value = 3618.4734
$ws.Range("M107").Value = -774.75
$ws.Range("N107").Value = -7458.4734

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 13892333
$ws.Range("I14").Value = 17858670
$ws.Range("K14").Value = 17858670
$ws.Range("M14").Value = -17858502
$ws.Range("H22").Value = 4127
$ws.Range("I22").Value = 5249.5
$ws.Range("J22").Value = 3004.5
$ws.Range("K22").Value = 5249.5
$ws.Range("L22").Value = 3004.5
$ws.Range("M22").Value = -4720.5
$ws.Range("N22").Value = -4062.5
$ws.Range("H41").Value = 7800
$ws.Range("I41").Value = 4400
$ws.Range("J41").Value = 9500
$ws.Range("K41").Value = 4400
$ws.Range("L41").Value = 9500
$ws.Range("M41").Value = -4045
$ws.Range("N41").Value = -10210
$ws.Range("H70").Value = 9685.571
$ws.Range("I70").Value = 5499.5
$ws.Range("J70").Value = 11360
$ws.Range("K70").Value = 5499.5
$ws.Range("L70").Value = 11360
$ws.Range("M70").Value = -5229.5
$ws.Range("N70").Value = -11900
$ws.Range("H73").Value = 9685.571
$ws.Range("I73").Value = 5499.5
$ws.Range("J73").Value = 11360
$ws.Range("K73").Value = 5499.5
$ws.Range("L73").Value = 11360
$ws.Range("M73").Value = -4563.5
$ws.Range("N73").Value = -13232
$ws.Range("H107").Value = 857.2308
$ws.Range("I107").Value = 433.57144
$ws.Range("K107").Value = 433.57144
$ws.Range("M107").Value = 1486.42856
$ws.Range("H113").Value = 3082.2144
$ws.Range("I113").Value = 3042.3845
$ws.Range("K113").Value = 3042.3845
$ws.Range("M113").Value = -872.3845000000001
$ws.Range("H122").Value = 2138.5
$ws.Range("I122").Value = 2366.375
$ws.Range("K122").Value = 7099.125
$ws.Range("M122").Value = -4649.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H40").Value = 3548.682
$ws.Range("I40").Value = 3176.4443
$ws.Range("K40").Value = 3176.4443
$ws.Range("M40").Value = -3040.4443
$ws.Range("H93").Value = 1998
$ws.Range("I93").Value = 1664
$ws.Range("K93").Value = 1664
$ws.Range("M93").Value = -416
$ws.Range("H122").Value = 3555.795
$ws.Range("I122").Value = 3253.25
$ws.Range("J122").Value = 4039.8667
$ws.Range("K122").Value = 9759.75
$ws.Range("L122").Value = 12119.6001
$ws.Range("M122").Value = -7309.75
$ws.Range("N122").Value = -17019.6001
$ws.Range("H136").Value = 5660.6523
$ws.Range("I136").Value = 4751.1333
$ws.Range("K136").Value = 14253.3999
$ws.Range("M136").Value = -11703.3999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 39998
$ws.Range("I37").Value = 39998
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 39998
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -39795
$ws.Range("N37").ClearContents()
$ws.Range("H122").Value = 7499.5
$ws.Range("I122").Value = 7499.5
$ws.Range("K122").Value = 22498.5
$ws.Range("M122").Value = -20048.5

